# Insert a new data row at row 710 (shifting existing rows 710-751 down to 711-752)
# and populate it with the new day's data: 2026/01/23, 金, 17, 138
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(710).Insert()

# Column A holds dates stored as plain text (e.g. "2026/01/23"), not real Excel
# date serials, so force text formatting before assigning the value to avoid
# automatic date conversion, then restore the default "Normal" style so the
# cell matches the unstyled look of its neighboring data cells.
$ws.Range("A710").NumberFormat = "@"
$ws.Range("A710").Value = "2026/01/23"
$ws.Range("A710").Style = "Normal"

$ws.Range("B710").Value = "金"
$ws.Range("C710").Value = 17
$ws.Range("D710").Value = 138
